$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.934.29"
$ws.Range("E2").Value = "  +2.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.379.64"
$ws.Range("E3").Value = "  +2.39%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.79"
$ws.Range("E5").Value = "  +2.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.20"
$ws.Range("E6").Value = "  +1.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +1.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.367.80"
$ws.Range("E8").Value = "  +2.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +12.21%  "

$ws.Range("E11").Value = "  +3.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.11"
$ws.Range("E12").Value = "  +2.84%  "

$ws.Range("E13").Value = "  +6.16%  "

$ws.Range("E14").Value = "  +3.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.920.49"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.21"
$ws.Range("E16").Value = "  +1.73%  "

$ws.Range("E17").Value = "  +2.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.379.73"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.918.85"
$ws.Range("E19").Value = "  +2.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.80"
$ws.Range("E20").Value = "  +1.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.991"
$ws.Range("E21").Value = "  +2.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.23"
$ws.Range("E22").Value = "  +14.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.94"
$ws.Range("E23").Value = "  +13.37%  "

$ws.Range("E24").Value = "  +2.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.84"
$ws.Range("E25").Value = "  +5.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.55"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("E27").Value = "  +7.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.75"
$ws.Range("E28").Value = "  +2.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.75"
$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.77"
$ws.Range("E30").Value = "  +6.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").Value = "  +5.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.50"
$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "570.08"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.32"
$ws.Range("E34").Value = "  +6.71%  "

$ws.Range("E35").Value = "  +2.41%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.66"
$ws.Range("E37").Value = "  +8.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").Value = "  -4.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.61"
$ws.Range("E39").Value = "  +2.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0748"
$ws.Range("E40").Value = "  +2.43%  "

$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.085.70"
$ws.Range("E42").Value = "  -0.68%  "

$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("E44").Value = "  +3.69%  "

$ws.Range("E45").Value = "  +4.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.134"
$ws.Range("E46").Value = "  +5.62%  "

$ws.Range("E47").Value = "  +2.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.14"
$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.20"
$ws.Range("E50").Value = "  +5.44%  "

$ws.Range("E51").Value = "  +3.79%  "
